# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 46075 (2026-02-22) to 46076 (2026-02-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C101").Value = 46076
